$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

# New orch_session_issue_id UUID values for rows 2-80 (column L)
$newIds = @(
    "ea4b34aa-1a69-4e1e-bf29-6f3b529abe00",
    "c951ee4c-8bd5-4d63-8349-76bd4baac952",
    "71a63c7c-0c0f-46be-87b9-926d0d008783",
    "d8be4f17-a53c-4ca2-a597-14ab2abd93da",
    "af9ac8bf-8c6e-4d63-b844-98b3047e6671",
    "45d47ab9-4026-4e9f-ba3a-4a45e1a6d1bd",
    "eeefb1f7-f9ee-41a3-8bb7-fddc7a101506",
    "4aaaae80-62ed-48e4-be19-b5d0cf220c83",
    "2d894ce0-a1f1-44e5-b711-43dc176b6d6c",
    "862a6e77-1083-4f88-904f-f4251fe37fb3",
    "4f9427dc-a8dc-46b4-86f6-27389849dd3a",
    "b452e9f0-d693-4434-922f-07610b0d7343",
    "60aba70d-c7fa-4bc0-8081-cfa13de10583",
    "772a619c-651f-463f-ba88-532452138666",
    "cbe85f61-a7c2-443f-ae3d-51f2acccdc12",
    "ab6c61b8-2a1d-4d14-bab7-75ecac746799",
    "410c3569-86a9-4bbc-927b-5bf07a20df1b",
    "c4799473-e0c4-4ad3-8d88-b57d53d1e2ca",
    "30a8e7e9-60ee-4d32-970f-a4adb4cabb33",
    "c005d20d-ef9b-4c75-a84e-e5e7228989c7",
    "bff980f0-a6d5-40c1-af1c-80e97b3fd67b",
    "9046b07a-1ab0-43e2-9eb2-9492c9591b24",
    "05e6cf6e-dd16-4153-9e2e-8bdbc4d5e466",
    "7eb80837-871f-43a1-aa5f-24da2735040c",
    "df824644-b31c-448b-b3fa-9ce906f90aa6",
    "7160b157-5599-4007-908c-1ac3066a452f",
    "45e479a7-0809-4d88-a279-31cba9128289",
    "daf9fc0f-c97e-4b23-9941-294d4931933c",
    "56d78d61-99c6-4627-b03c-cf070b31dbfd",
    "214f5c1d-90d5-43cf-b48e-92694262fd89",
    "4f4fb0ec-d17c-4756-8889-f2a27eae3a28",
    "f2e98b1a-97e7-4a56-bdbb-cf024bac813f",
    "f573016b-c933-4694-a254-fdd975460a59",
    "1fa35da9-c562-4418-ab5c-9675bcb0c169",
    "5b4721d1-5ebe-4dcc-a4e4-11f569607679",
    "2bab4fd3-893a-4cae-a47a-03a5f52e81e3",
    "7e10905c-a5f2-4d8d-816c-46939bc7f288",
    "1875bacf-0fea-4743-b355-20f897c62328",
    "398367a3-8e6f-4db1-8e4a-a27685203c79",
    "1c99b50d-3a59-4ab6-b58e-ad21188c3835",
    "fe351f07-0fa9-43ce-94d0-2f69ef6a9209",
    "c53bfa9f-b87c-4a90-a7e4-dfcab376428e",
    "61cf0f84-634e-46f2-917c-cf99150fe9ab",
    "47667083-d500-42f8-95e1-f7b14686b5b6",
    "eb795860-b156-42f0-84eb-e4ec4e606efb",
    "3210278e-c460-4bbc-bbcb-6a1b74a37b74",
    "0f64d539-1d71-4ae6-a730-f08cc7349eef",
    "dbac4637-04f7-4f73-967b-8e3aa0a55718",
    "cd9cf62e-67d5-4245-81c4-0080f07e39eb",
    "1c8fdd4f-a1f7-4a19-a824-ba621a45b257",
    "51caecb9-122d-4ba1-902d-e9e047caa0f7",
    "db8a38ed-cf83-400e-a4a5-192ea849872f",
    "b9bee7a3-f016-43ab-9b6b-f35bc552c3c2",
    "1f5d6c79-83c9-47ca-bf98-d76b4b6ead4a",
    "7eb77061-590d-429c-9b4b-64b55cfe4dd4",
    "18c10dd8-fd51-442a-9dbb-39d1c0167fbf",
    "186ca05a-075a-4794-add0-3213127a815c",
    "1223e4d2-10d1-4dff-b7eb-e3827f0a8f9d",
    "fc437c56-50e1-40fa-aa1f-80e79639e079",
    "1cd94cbe-b70e-4116-9df1-38fc42f61b0b",
    "23d8e166-7251-4f8f-9318-834ae37de8f8",
    "a7ca19f9-a69a-49bf-8f59-74ee70c52ee5",
    "7e8dab49-1a1e-4129-9a02-a8792e278975",
    "a0d0cde7-405d-4fb8-9dc6-cca1c554c20e",
    "52a68fd9-0fe1-47d1-865f-7bf9a95adb3f",
    "2069624f-6b41-43f1-b4c1-1bac0a8a095d",
    "95d07939-143d-4313-bbac-884627573ef7",
    "f47f4b43-2b14-4f01-b7a6-8be0fe0a7e3b",
    "618e3bf6-8b2c-44b7-bf51-9695ae08356b",
    "b7bc73cc-e366-4c33-b070-5ad394bcd805",
    "ff083d1b-81e3-4ae4-97b7-1f95df6b650a",
    "8b475a72-35ee-4570-b1c1-ca748fdd1ba0",
    "dd17a62a-c8dd-4450-840c-c870137c0916",
    "51521b14-8912-44e5-acf5-16a9c3e29ecd",
    "96615a90-e365-4c1c-ac13-f47a36c7f607",
    "2257fd37-edaf-4d07-b823-5a0174f51bd5",
    "ad9c42f7-ec3f-4310-8e61-33e5d29716eb",
    "0fc17170-2adb-4206-af1c-f8c8f04c81cf",
    "a117ea99-445c-4a50-a703-081afc88f923"
)

for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $newIds[$i]
}

# New orch_session_issue_id UUID value for row 96 (column L)
$ws.Cells.Item(96, 12).Value = "402e96de-6b86-4346-a47b-259fb06965bf"

